# Updates the Coinranking crypto price/volume snapshot in-place to match the
# latest scrape (GitHub Actions refresh). Rows 2-39 get refreshed Price/
# Volume(1h) figures; rows 40-51 shift to the next batch of coins (Frax drops
# off the bottom, EOS is newly appended).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "27.930.07"
$ws.Range("D3").Value = "1.882.61"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("D5").Value = "'335.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  +1.30%  "
$ws.Range("D7").Value = "'0.4680"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("D8").Value = "'0.3915"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").Value = "'47.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").Value = "'0.07948"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "'1.011"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").Value = "'21.65"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "1.894.75"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("D14").Value = "'5.952"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").Value = "'7.110"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "'1.020"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").Value = "'0.06805"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").Value = "'87.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "'0.00001045"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("D21").Value = "'1.017"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").Value = "27.943.92"
$ws.Range("E22").Value = "  +0.69%  "
$ws.Range("D23").Value = "'5.471"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("D25").Value = "'2.355"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("D26").Value = "2.129.83"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").Value = "'159.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").Value = "'20.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'2.071"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "'5.448"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("D31").Value = "'120.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").Value = "'0.09537"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").Value = "'0.9566"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("D34").Value = "'3.657"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Value = "'5.317"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'1.349"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.37%  "
$ws.Range("D37").Value = "'0.06118"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "'0.02232"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Value = "'1.204"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'8.121"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5867"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1892"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'10.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.20%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.272"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.64%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5626"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.401"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.919"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06857"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'113.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "'1.065"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.89%  "
